# Regenerate the s_vals data (filter save games) -- update B:E and G for
# rows 2-20 on the active sheet with the freshly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = 3.272327238179451;    C = 1.626987699542094;    D = 0.7210945179870265;  E = 0.5333859586016987;  G = 6.15379541431027 }
    @{ Row = 3;  B = 3.272327238179451;    C = 1.626987699542094;    D = 0.1496068669990043;   E = 0.5333859586016987;  G = 5.582307763322248 }
    @{ Row = 4;  B = 0.6545652718822623;   C = 1.626987699542094;    D = 3.223369029078222;    E = 0.5333859586016987;  G = 6.038307959104277 }
    @{ Row = 5;  B = 3.272327238179451;    C = 1.626987699542094;    D = 3.223369029078222;    E = 0.5333859586016987;  G = 8.656069925401464 }
    @{ Row = 6;  B = 0.003078177322033415; C = 0.00006708468553440206; D = 0.1496068669990043;  E = 0.5333859586016987;  G = 0.6861380876082708 }
    @{ Row = 7;  B = 1.445647641019636;    C = 1.626987699542094;    D = 3.223369029078222;    E = 13.86384647080068;   G = 20.15985084044064 }
    @{ Row = 8;  B = 1.445647641019636;    C = 1.626987699542094;    D = 0.1496068669990043;   E = 0.5333859586016987;  G = 3.755628166162433 }
    @{ Row = 9;  B = 0.2881169905109251;   C = 0.04103571897497393;  D = 0.7210945179870265;   E = 0.5333859586016987;  G = 1.583633186074624 }
    @{ Row = 10; B = 0.04172184405617529;  C = 0.3048912486333797;   D = 0.1496068669990043;   E = 0.5333859586016987;  G = 1.029605918290258 }
    @{ Row = 11; B = 0.01253208636536152;  C = 0.002658071450198252; D = 0.1496068669990043;   E = 0.5333859586016987;  G = 0.6981829834162627 }
    @{ Row = 12; B = 3.272327238179451;    C = 1.626987699542094;    D = 0.7210945179870265;   E = 0.5333859586016987;  G = 6.15379541431027 }
    @{ Row = 13; B = 3.272327238179451;    C = 1.626987699542094;    D = 0.7210945179870265;   E = 0.5333859586016987;  G = 6.15379541431027 }
    @{ Row = 14; B = 3.272327238179451;    C = 1.626987699542094;    D = 3.223369029078222;    E = 0.5333859586016987;  G = 8.656069925401464 }
    @{ Row = 15; B = 0.6545652718822623;   C = 1.626987699542094;    D = 3.223369029078222;    E = 13.86384647080068;   G = 19.36876847130326 }
    @{ Row = 16; B = 0.6545652718822623;   C = 1.626987699542094;    D = 0.1496068669990043;   E = 0.5333859586016987;  G = 2.964545797025059 }
    @{ Row = 17; B = 0.1169995834814548;   C = 0.3048912486333797;   D = 0.1496068669990043;   E = 0.5333859586016987;  G = 1.104883657715537 }
    @{ Row = 18; B = 0.6545652718822623;   C = 0.3048912486333797;   D = 0.7210945179870265;   E = 0.5333859586016987;  G = 2.213936997104367 }
    @{ Row = 19; B = 0.6545652718822623;   C = 1.626987699542094;    D = 0.7210945179870265;   E = 0.5333859586016987;  G = 3.536033448013082 }
    @{ Row = 20; B = 0.6545652718822623;   C = 1.626987699542094;    D = 3.223369029078222;    E = 0.5333859586016987;  G = 6.038307959104277 }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 7).Value = $row.G
}
